$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "P4 - Roy"

# --- Enter the new shared-text values in the same order the original author
# typed them, so the regenerated sharedStrings table lines up with the target. ---
$ws.Range("A50").Value = "Export PDF "
$ws.Range("A52").Value = "PDF document css"
$ws.Range("A53").Value = "PDF document css"
$ws.Range("D50").Value = "afronden pdf en begin css"
$ws.Range("A51").Value = "Css leren"
$ws.Range("D51").Value = "online informatie opzoeken en begin maken"
$ws.Range("D52").Value = "zelf proberen te stylen"
$ws.Range("D54").Value = "Pdf document afgerond en gestyled"
$ws.Range("D53").Value = "Hulp gekregen van Ivar met wat start problemen"
$ws.Range("A54").Value = "Afronden PDF document"

# --- Dates: B50 gets a brand-new "d-mmm" number format (style index 13);
# B51:B54 reuse the existing date style already used by B47:B49 (style index 3). ---
$ws.Range("B50").NumberFormat = "d-mmm"
$ws.Range("B50").Value = 44163   # 28-Nov-2020

$ws.Range("B47").Copy()
$ws.Range("B51").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B52").PasteSpecial(-4122)
$ws.Range("B53").PasteSpecial(-4122)
$ws.Range("B54").PasteSpecial(-4122)

$ws.Range("B51").Value = 44164   # 29-Nov-2020
$ws.Range("B52").Value = 44165   # 30-Nov-2020
$ws.Range("B53").Value = 44165
$ws.Range("B54").Value = 44165

# --- Time spent (minutes) ---
$ws.Range("C50").Value = 90
$ws.Range("C51").Value = 90
$ws.Range("C52").Value = 60
$ws.Range("C53").Value = 40
$ws.Range("C54").Value = 150

# --- Grow the worksheet table to cover the 5 new rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A9:D54"))

# --- Refresh the hour-total formulas so they include the new rows ---
$ws.Range("B6").Formula = "=SUM(C10:C154)/60"
$ws.Range("B7").Formula = "=(2*3*28)-(SUM(C10:C154)/60)"
$excel.CalculateFull()

# --- Make this the active/visible sheet, like in the saved workbook ---
$ws.Activate()
$ws.Range("C55").Select()
